$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order Data")

# Row 2 (Order ID 1): mark completed, assign waiter 7, set status to Completed
$ws.Range("D2").Value = $true
$ws.Range("E2").Value = "Completed"
$ws.Range("G2").Value = 7

# Row 6 (Order ID 5): mark completed, set status to Completed
$ws.Range("D6").Value = $true
$ws.Range("E6").Value = "Completed"

[void]$ws.Activate()
[void]$ws.Range("E9").Select()
